# "fix: relatorio dash kpi"
#
# The Dashboard sheet's A1 cell used to hold a leftover long-date
# number format with no content. Turn it into the KPI header's
# "report generated at" label, and give the adjacent B1 cell a
# short-date format to hold the generation date, moving the active
# selection there.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Dashboard")

# A1: drop the stale formatting and add the label text.
$a1 = $ws.Range("A1")
$a1.ClearFormats()
$a1.Value2 = "Relatório gerado em:"

# B1: new cell, formatted as a short date (left without a value).
$b1 = $ws.Range("B1")
$b1.NumberFormat = "mm-dd-yy"

# The active selection moves from B4 to B1.
$b1.Select()
